{"js": "// Title paragraph \"Hello World Lab\" -> \"Hello World\": drop the trailing\n// \" Lab\" (removes the \" \" and \"Lab\" runs, matching the XML diff).\nconst titleHits = context.document.body.search(\" Lab\", { matchCase: true, matchWholeWord: false });\ntitleHits.load(\"items\");\nawait context.sync();\n\nif (titleHits.items.length > 0) {\n  titleHits.items[0].delete();\n  await context.sync();\n}\n\n// Date paragraph: update the printed time from 08:45:09 PM to 08:54:51 PM.\nconst dateHits = context.document.body.search(\"August  10, 2021 (08:45:09 PM)\", { matchCase: true, matchWholeWord: false });\ndateHits.load(\"items\");\nawait context.sync();\n\nif (dateHits.items.length > 0) {\n  dateHits.items[0].insertText(\"August  10, 2021 (08:54:51 PM)\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Title paragraph \"Hello World Lab\" -> \"Hello World\": drop the trailing\n# \" Lab\" text (removes the \" \" and \"Lab\" runs, matching the XML diff).\n# MatchCase avoids accidentally matching other case-insensitive \"lab\"\n# occurrences elsewhere in the document body.\n$titleRange = $d.Content\n$titleRange.Find.ClearFormatting()\n$titleRange.Find.Text = \" Lab\"\n$titleRange.Find.MatchCase = $true\n$titleRange.Find.MatchWholeWord = $false\n$titleRange.Find.Forward = $true\n$foundTitle = $titleRange.Find.Execute()\nif ($foundTitle) {\n    $titleRange.Text = \"\"\n}\n\n# Date paragraph: update the printed time from 08:45:09 PM to 08:54:51 PM.\n# Replace the whole paragraph's Range (paragraph mark included) so the\n# surviving run keeps its xml:space=\"preserve\" text run, same as the source.\n$datePara = $d.Paragraphs.Item(3)\n$dateRange = $datePara.Range\nif ($dateRange.Text -match \"August  10, 2021 \\(08:45:09 PM\\)\") {\n    $dateRange.Text = \"August  10, 2021 (08:54:51 PM)\"\n}\n"}
